# Apply leaderboard update to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("G7").Value = 4.0
$ws.Range("H7").Value = 4.0
$ws.Range("J7").Value = "Yes"

# Row 8
$ws.Range("G8").Value = 1.0

# Row 9
$ws.Range("G9").Value = 3.0

# Row 11
$ws.Range("G11").Value = 4.0
$ws.Range("H11").Value = 4.0
$ws.Range("J11").Value = "Yes"

# Row 16
$ws.Range("G16").Value = 4.0
$ws.Range("H16").Value = 4.0
$ws.Range("J16").Value = "Yes"

# Row 24
$ws.Range("G24").Value = 2.0

# Row 47
$ws.Range("G47").Value = 1.0

# Row 49
$ws.Range("G49").Value = 1.0

# Row 64
$ws.Range("G64").Value = 1.0

# Row 69
$ws.Range("G69").Value = 1.0

# Row 78
$ws.Range("G78").Value = 2.0

# Row 85
$ws.Range("K85").Value = "Yes"

# Row 95
$ws.Range("G95").Value = 1.0

# Row 103
$ws.Range("G103").Value = 4.0
$ws.Range("H103").Value = 1.0

# Row 110
$ws.Range("G110").Value = 4.0
$ws.Range("H110").Value = 1.0

# Row 122
$ws.Range("H122").Value = 1.0

# Row 137
$ws.Range("K137").Value = "Yes"

# Row 145
$ws.Range("H145").Value = 4.0
$ws.Range("J145").Value = "Yes"

# Row 151
$ws.Range("G151").Value = 1.0
$ws.Range("H151").Value = 2.0

# Row 179
$ws.Range("G179").Value = 1.0
